# Applies the scheduled-runner price/profit refresh to the Anima_Profits workbook.
# Each of the 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) shares the same
# Table_<Job> layout: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ, N=LeveProfitHQ.
# Update a handful of rows per sheet with refreshed market-board figures;
# ClearContents() is used where a column no longer has a value (matches source diff,
# which drops the <c> element entirely rather than writing an empty/zero cell).
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item(1)
# Row 86
$ws.Cells.Item(86, 8).Value = 220067400
$ws.Cells.Item(86, 9).Value = 220067400
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 220067400
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 13).ClearContents()
$ws.Cells.Item(86, 14).Value = -220066277

# Row 89
$ws.Cells.Item(89, 8).Value = 220067400
$ws.Cells.Item(89, 9).Value = 220067400
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 11).Value = 1100337000
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 13).ClearContents()
$ws.Cells.Item(89, 14).Value = -1100331384

# Row 112
$ws.Cells.Item(112, 8).Value = 5084.3413
$ws.Cells.Item(112, 10).Value = 5589.1353
$ws.Cells.Item(112, 12).Value = 16767.4059
$ws.Cells.Item(112, 14).Value = -18983.4059

# Row 113
$ws.Cells.Item(113, 8).Value = 2447.4
$ws.Cells.Item(113, 9).Value = 2238.125
$ws.Cells.Item(113, 10).Value = 2686.5715
$ws.Cells.Item(113, 11).Value = 2238.125
$ws.Cells.Item(113, 12).Value = 2686.5715
$ws.Cells.Item(113, 13).Value = 1015.875
$ws.Cells.Item(113, 14).Value = -9194.5715

# Row 132
$ws.Cells.Item(132, 8).Value = 1428.8654
$ws.Cells.Item(132, 9).Value = 1256.0454
$ws.Cells.Item(132, 10).Value = 2379.375
$ws.Cells.Item(132, 11).Value = 3768.1362
$ws.Cells.Item(132, 12).Value = 7138.125
$ws.Cells.Item(132, 13).Value = -1238.1362
$ws.Cells.Item(132, 14).Value = -12198.125

# Row 137
$ws.Cells.Item(137, 8).Value = 1812863.8
$ws.Cells.Item(137, 9).Value = 3788766
$ws.Cells.Item(137, 10).Value = 1620.1666
$ws.Cells.Item(137, 11).Value = 11366298
$ws.Cells.Item(137, 12).Value = 4860.4998
$ws.Cells.Item(137, 13).Value = -11363748
$ws.Cells.Item(137, 14).Value = -9960.4998

# Row 140
$ws.Cells.Item(140, 8).Value = 72182.46000000001
$ws.Cells.Item(140, 10).Value = 72182.46000000001
$ws.Cells.Item(140, 12).Value = 72182.46000000001
$ws.Cells.Item(140, 14).Value = -82542.46000000001

# Row 141
$ws.Cells.Item(141, 8).Value = 3493.96
$ws.Cells.Item(141, 9).Value = 1681.3334
$ws.Cells.Item(141, 10).Value = 8155
$ws.Cells.Item(141, 11).Value = 5044.0002
$ws.Cells.Item(141, 12).Value = 24465
$ws.Cells.Item(141, 13).Value = 135.9997999999996
$ws.Cells.Item(141, 14).Value = -34825

# --- ARM sheet ---
$ws = $wb.Worksheets.Item(2)
# Row 4
$ws.Cells.Item(4, 8).Value = 0
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 12).ClearContents()
$ws.Cells.Item(4, 14).Value = 0

# Row 5
$ws.Cells.Item(5, 8).Value = 192.2
$ws.Cells.Item(5, 9).Value = 101
$ws.Cells.Item(5, 10).Value = 215
$ws.Cells.Item(5, 11).Value = 101
$ws.Cells.Item(5, 12).Value = 215
$ws.Cells.Item(5, 13).Value = 11
$ws.Cells.Item(5, 14).Value = -439

# Row 32
$ws.Cells.Item(32, 8).Value = 6680387.5
$ws.Cells.Item(32, 9).Value = 7824461
$ws.Cells.Item(32, 11).Value = 7824461
$ws.Cells.Item(32, 13).Value = -7824174

# Row 45
$ws.Cells.Item(45, 8).Value = 3047.0833
$ws.Cells.Item(45, 9).Value = 2876
$ws.Cells.Item(45, 11).Value = 2876
$ws.Cells.Item(45, 13).Value = -2499

# Row 74
$ws.Cells.Item(74, 8).Value = 13160367
$ws.Cells.Item(74, 9).Value = 1748.15
$ws.Cells.Item(74, 10).Value = 27781056
$ws.Cells.Item(74, 11).Value = 1748.15
$ws.Cells.Item(74, 12).Value = 27781056
$ws.Cells.Item(74, 13).Value = -874.1500000000001
$ws.Cells.Item(74, 14).Value = -27782804

# Row 77
$ws.Cells.Item(77, 8).Value = 13160367
$ws.Cells.Item(77, 9).Value = 1748.15
$ws.Cells.Item(77, 10).Value = 27781056
$ws.Cells.Item(77, 11).Value = 8740.75
$ws.Cells.Item(77, 12).Value = 138905280
$ws.Cells.Item(77, 13).Value = -4372.75
$ws.Cells.Item(77, 14).Value = -138914016

# Row 106
$ws.Cells.Item(106, 8).Value = 74185
$ws.Cells.Item(106, 10).Value = 74185
$ws.Cells.Item(106, 12).Value = 74185
$ws.Cells.Item(106, 14).Value = -76709

# Row 140
$ws.Cells.Item(140, 8).Value = 76795.5
$ws.Cells.Item(140, 10).Value = 76795.5
$ws.Cells.Item(140, 12).Value = 76795.5
$ws.Cells.Item(140, 14).Value = -87155.5

# --- BSM sheet ---
$ws = $wb.Worksheets.Item(3)
# Row 4
$ws.Cells.Item(4, 8).Value = 192.2
$ws.Cells.Item(4, 9).Value = 101
$ws.Cells.Item(4, 10).Value = 215
$ws.Cells.Item(4, 11).Value = 101
$ws.Cells.Item(4, 12).Value = 215
$ws.Cells.Item(4, 13).Value = 14
$ws.Cells.Item(4, 14).Value = -445

# Row 22
$ws.Cells.Item(22, 8).Value = 765
$ws.Cells.Item(22, 9).Value = 818
$ws.Cells.Item(22, 10).Value = 500
$ws.Cells.Item(22, 11).Value = 818
$ws.Cells.Item(22, 12).Value = 500
$ws.Cells.Item(22, 13).Value = -645
$ws.Cells.Item(22, 14).Value = -846

# Row 86
$ws.Cells.Item(86, 8).Value = 1936.0416
$ws.Cells.Item(86, 9).Value = 1930.2273
$ws.Cells.Item(86, 10).Value = 2000
$ws.Cells.Item(86, 11).Value = 1930.2273
$ws.Cells.Item(86, 12).Value = 2000
$ws.Cells.Item(86, 13).Value = -807.2273
$ws.Cells.Item(86, 14).Value = -4246

# Row 89
$ws.Cells.Item(89, 8).Value = 1936.0416
$ws.Cells.Item(89, 9).Value = 1930.2273
$ws.Cells.Item(89, 10).Value = 2000
$ws.Cells.Item(89, 11).Value = 9651.136500000001
$ws.Cells.Item(89, 12).Value = 10000
$ws.Cells.Item(89, 13).Value = -4035.136500000001
$ws.Cells.Item(89, 14).Value = -21232

# Row 134
$ws.Cells.Item(134, 8).Value = 2735.3142
$ws.Cells.Item(134, 9).Value = 2055.7273
$ws.Cells.Item(134, 10).Value = 3885.3845
$ws.Cells.Item(134, 11).Value = 6167.1819
$ws.Cells.Item(134, 12).Value = 11656.1535
$ws.Cells.Item(134, 13).Value = -3632.1819
$ws.Cells.Item(134, 14).Value = -16726.1535

# --- CRP sheet ---
$ws = $wb.Worksheets.Item(4)
# Row 7
$ws.Cells.Item(7, 8).Value = 72.72727
$ws.Cells.Item(7, 9).Value = 78.2
$ws.Cells.Item(7, 10).Value = 68.166664
$ws.Cells.Item(7, 11).Value = 78.2
$ws.Cells.Item(7, 12).Value = 68.166664
$ws.Cells.Item(7, 13).Value = 34.8
$ws.Cells.Item(7, 14).Value = -294.166664

# Row 16
$ws.Cells.Item(16, 8).Value = 1642.1111
$ws.Cells.Item(16, 9).Value = 1520.9
$ws.Cells.Item(16, 11).Value = 1520.9
$ws.Cells.Item(16, 13).Value = -1233.9

# Row 31
$ws.Cells.Item(31, 8).Value = 5259.039
$ws.Cells.Item(31, 9).Value = 1705.8695
$ws.Cells.Item(31, 10).Value = 6772.426
$ws.Cells.Item(31, 11).Value = 1705.8695
$ws.Cells.Item(31, 12).Value = 6772.426
$ws.Cells.Item(31, 13).Value = -1410.8695
$ws.Cells.Item(31, 14).Value = -7362.426

# Row 34
$ws.Cells.Item(34, 8).Value = 5259.039
$ws.Cells.Item(34, 9).Value = 1705.8695
$ws.Cells.Item(34, 10).Value = 6772.426
$ws.Cells.Item(34, 11).Value = 1705.8695
$ws.Cells.Item(34, 12).Value = 6772.426
$ws.Cells.Item(34, 13).Value = -1503.8695
$ws.Cells.Item(34, 14).Value = -7176.426

# Row 113
$ws.Cells.Item(113, 8).Value = 1642.1111
$ws.Cells.Item(113, 9).Value = 1520.9
$ws.Cells.Item(113, 11).Value = 1520.9
$ws.Cells.Item(113, 13).Value = 649.0999999999999

# Row 132
$ws.Cells.Item(132, 8).Value = 10640049
$ws.Cells.Item(132, 9).Value = 13890294
$ws.Cells.Item(132, 10).Value = 2882.818
$ws.Cells.Item(132, 11).Value = 41670882
$ws.Cells.Item(132, 12).Value = 8648.454000000002
$ws.Cells.Item(132, 13).Value = -41668352
$ws.Cells.Item(132, 14).Value = -13708.454

# --- CUL sheet ---
$ws = $wb.Worksheets.Item(5)
# Row 5
$ws.Cells.Item(5, 8).Value = 875.81134
$ws.Cells.Item(5, 10).Value = 1214.7646
$ws.Cells.Item(5, 12).Value = 3644.2938
$ws.Cells.Item(5, 14).Value = -3868.2938

# Row 131
$ws.Cells.Item(131, 8).Value = 3926.718
$ws.Cells.Item(131, 10).Value = 5270.5713
$ws.Cells.Item(131, 12).Value = 15811.7139
$ws.Cells.Item(131, 14).Value = -25891.7139

# Row 135
$ws.Cells.Item(135, 8).Value = 875.81134
$ws.Cells.Item(135, 10).Value = 1214.7646
$ws.Cells.Item(135, 12).Value = 10932.8814
$ws.Cells.Item(135, 14).Value = -16002.8814

# Row 141
$ws.Cells.Item(141, 8).Value = 7734.7827
$ws.Cells.Item(141, 9).Value = 3878.889
$ws.Cells.Item(141, 10).Value = 10213.571
$ws.Cells.Item(141, 11).Value = 11636.667
$ws.Cells.Item(141, 12).Value = 30640.713
$ws.Cells.Item(141, 13).Value = -6456.667000000001
$ws.Cells.Item(141, 14).Value = -41000.713

# --- GSM sheet ---
$ws = $wb.Worksheets.Item(6)
# Row 2
$ws.Cells.Item(2, 8).Value = 147.6
$ws.Cells.Item(2, 9).Value = 202.77777
$ws.Cells.Item(2, 10).Value = 64.833336
$ws.Cells.Item(2, 11).Value = 202.77777
$ws.Cells.Item(2, 12).Value = 64.833336
$ws.Cells.Item(2, 13).Value = -89.77777
$ws.Cells.Item(2, 14).Value = -290.833336

# --- LTW sheet ---
$ws = $wb.Worksheets.Item(7)
# Row 61
$ws.Cells.Item(61, 8).Value = 2772.76
$ws.Cells.Item(61, 9).Value = 1201.1875
$ws.Cells.Item(61, 10).Value = 5566.6665
$ws.Cells.Item(61, 11).Value = 1201.1875
$ws.Cells.Item(61, 12).Value = 5566.6665
$ws.Cells.Item(61, 13).Value = -999.1875
$ws.Cells.Item(61, 14).Value = -5970.6665

# Row 113
$ws.Cells.Item(113, 8).Value = 2772.76
$ws.Cells.Item(113, 9).Value = 1201.1875
$ws.Cells.Item(113, 10).Value = 5566.6665
$ws.Cells.Item(113, 11).Value = 1201.1875
$ws.Cells.Item(113, 12).Value = 5566.6665
$ws.Cells.Item(113, 13).Value = 968.8125
$ws.Cells.Item(113, 14).Value = -9906.666499999999

# Row 136
$ws.Cells.Item(136, 8).Value = 1927.2778
$ws.Cells.Item(136, 9).Value = 1691.2307
$ws.Cells.Item(136, 10).Value = 2541
$ws.Cells.Item(136, 11).Value = 5073.6921
$ws.Cells.Item(136, 12).Value = 7623
$ws.Cells.Item(136, 13).Value = -2523.6921
$ws.Cells.Item(136, 14).Value = -12723

# Row 140
$ws.Cells.Item(140, 8).Value = 57266.168
$ws.Cells.Item(140, 10).Value = 57266.168
$ws.Cells.Item(140, 12).Value = 57266.168
$ws.Cells.Item(140, 14).Value = -67626.16800000001

# --- WVR sheet ---
$ws = $wb.Worksheets.Item(8)
# Row 92
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).ClearContents()
$ws.Cells.Item(92, 14).Value = 0

# Row 93
$ws.Cells.Item(93, 8).Value = 78389
$ws.Cells.Item(93, 10).Value = 78389
$ws.Cells.Item(93, 12).Value = 78389
$ws.Cells.Item(93, 14).Value = -83381

# Row 105
$ws.Cells.Item(105, 8).Value = 94307.5
$ws.Cells.Item(105, 10).Value = 94307.5
$ws.Cells.Item(105, 12).Value = 94307.5
$ws.Cells.Item(105, 14).Value = -101295.5

# Row 122
$ws.Cells.Item(122, 8).Value = 1839.25
$ws.Cells.Item(122, 9).Value = 1793.9048
$ws.Cells.Item(122, 10).Value = 2156.6667
$ws.Cells.Item(122, 11).Value = 5381.7144
$ws.Cells.Item(122, 12).Value = 6470.000100000001
$ws.Cells.Item(122, 13).Value = -2931.7144
$ws.Cells.Item(122, 14).Value = -11370.0001

# Row 136
$ws.Cells.Item(136, 8).Value = 3119.9805
$ws.Cells.Item(136, 9).Value = 2749.5151
$ws.Cells.Item(136, 10).Value = 3799.1667
$ws.Cells.Item(136, 11).Value = 8248.5453
$ws.Cells.Item(136, 12).Value = 11397.5001
$ws.Cells.Item(136, 13).Value = -5698.5453
$ws.Cells.Item(136, 14).Value = -16497.5001

